$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.162.22'
$ws.Range("E2").Value = '  -3.28%  '
$ws.Range("D3").Value = '1.596.92'
$ws.Range("E3").Value = '  -3.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.75'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3766'
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3659'
$ws.Range("E8").Value = '  -4.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.72'
$ws.Range("E9").Value = '  -4.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.004'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.276'
$ws.Range("E11").Value = '  -5.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08092'
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.88'
$ws.Range("E13").Value = '  -4.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.636'
$ws.Range("E14").Value = '  -7.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.570'
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001261'
$ws.Range("E16").Value = '  -3.70%  '
$ws.Range("D17").Value = '1.595.50'
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.10'
$ws.Range("E18").Value = '  -2.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06814'
$ws.Range("E19").Value = '  -2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.48'
$ws.Range("E20").Value = '  -6.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.607'
$ws.Range("E21").Value = '  -4.28%  '
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.12'
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("D24").Value = '23.163.24'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.359'
$ws.Range("E25").Value = '  -5.08%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.958'
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.13'
$ws.Range("E27").Value = '  -4.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.03'
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.229'
$ws.Range("E29").Value = '  -4.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.43'
$ws.Range("E30").Value = '  -5.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.467'
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.108'
$ws.Range("E32").Value = '  -8.87%  '
$ws.Range("D33").Value = '1.771.90'
$ws.Range("E33").Value = '  -3.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9749'
$ws.Range("E34").Value = '  -5.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07748'
$ws.Range("E35").Value = '  -3.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02780'
$ws.Range("E36").Value = '  -6.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.260'
$ws.Range("E37").Value = '  -6.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2548'
$ws.Range("E38").Value = '  -5.24%  '
$ws.Range("E39").Value = '  -7.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08871'
$ws.Range("E40").Value = '  -2.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.385'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7165'
$ws.Range("E42").Value = '  -5.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.79'
$ws.Range("E43").Value = '  -5.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.12'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6637'
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.315'
$ws.Range("E46").Value = '  -5.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.966'
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07987'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.23'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.172'
$ws.Range("E51").Value = '  -4.87%  '
